{"js": "// Replace each old three-digit-division answer with its new value.\n// Source strings are each unique within the document, so a direct\n// body.search() + insertText(replace) round-trip is unambiguous and\n// preserves the existing run formatting (font/size) of each cell.\nconst replacements = [\n  [\"957\u00f76=159, 3\", \"413\u00f74=103, 1\"],\n  [\"974\u00f72=487, 0\", \"574\u00f79=63, 7\"],\n  [\"829\u00f72=414, 1\", \"619\u00f72=309, 1\"],\n  [\"975\u00f79=108, 3\", \"881\u00f79=97, 8\"],\n  [\"161\u00f77=23, 0\", \"542\u00f74=135, 2\"],\n  [\"780\u00f78=97, 4\", \"453\u00f72=226, 1\"],\n  [\"696\u00f77=99, 3\", \"583\u00f72=291, 1\"],\n  [\"880\u00f78=110, 0\", \"169\u00f77=24, 1\"],\n  [\"846\u00f72=423, 0\", \"640\u00f75=128, 0\"],\n  [\"676\u00f73=225, 1\", \"951\u00f73=317, 0\"],\n  [\"481\u00f78=60, 1\", \"596\u00f77=85, 1\"],\n  [\"611\u00f74=152, 3\", \"717\u00f73=239, 0\"],\n  [\"528\u00f78=66, 0\", \"155\u00f76=25, 5\"],\n  [\"664\u00f75=132, 4\", \"261\u00f73=87, 0\"],\n  [\"768\u00f77=109, 5\", \"519\u00f76=86, 3\"],\n  [\"259\u00f74=64, 3\", \"296\u00f72=148, 0\"],\n  [\"462\u00f74=115, 2\", \"442\u00f76=73, 4\"],\n  [\"786\u00f79=87, 3\", \"903\u00f73=301, 0\"],\n  [\"101\u00f72=50, 1\", \"826\u00f72=413, 0\"],\n  [\"850\u00f78=106, 2\", \"893\u00f72=446, 1\"],\n  [\"642\u00f72=321, 0\", \"414\u00f73=138, 0\"],\n  [\"928\u00f76=154, 4\", \"968\u00f73=322, 2\"],\n  [\"467\u00f73=155, 2\", \"690\u00f77=98, 4\"],\n  [\"203\u00f73=67, 2\", \"356\u00f77=50, 6\"],\n  [\"289\u00f78=36, 1\", \"793\u00f77=113, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-division answer with its new value.\n# Each source string is unique within the document, so Find/Execute with\n# Replace:=wdReplaceOne against the whole document Range is unambiguous\n# and leaves the existing run formatting (font/size) of each cell intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"957\u00f76=159, 3\", \"413\u00f74=103, 1\"),\n    @(\"974\u00f72=487, 0\", \"574\u00f79=63, 7\"),\n    @(\"829\u00f72=414, 1\", \"619\u00f72=309, 1\"),\n    @(\"975\u00f79=108, 3\", \"881\u00f79=97, 8\"),\n    @(\"161\u00f77=23, 0\", \"542\u00f74=135, 2\"),\n    @(\"780\u00f78=97, 4\", \"453\u00f72=226, 1\"),\n    @(\"696\u00f77=99, 3\", \"583\u00f72=291, 1\"),\n    @(\"880\u00f78=110, 0\", \"169\u00f77=24, 1\"),\n    @(\"846\u00f72=423, 0\", \"640\u00f75=128, 0\"),\n    @(\"676\u00f73=225, 1\", \"951\u00f73=317, 0\"),\n    @(\"481\u00f78=60, 1\", \"596\u00f77=85, 1\"),\n    @(\"611\u00f74=152, 3\", \"717\u00f73=239, 0\"),\n    @(\"528\u00f78=66, 0\", \"155\u00f76=25, 5\"),\n    @(\"664\u00f75=132, 4\", \"261\u00f73=87, 0\"),\n    @(\"768\u00f77=109, 5\", \"519\u00f76=86, 3\"),\n    @(\"259\u00f74=64, 3\", \"296\u00f72=148, 0\"),\n    @(\"462\u00f74=115, 2\", \"442\u00f76=73, 4\"),\n    @(\"786\u00f79=87, 3\", \"903\u00f73=301, 0\"),\n    @(\"101\u00f72=50, 1\", \"826\u00f72=413, 0\"),\n    @(\"850\u00f78=106, 2\", \"893\u00f72=446, 1\"),\n    @(\"642\u00f72=321, 0\", \"414\u00f73=138, 0\"),\n    @(\"928\u00f76=154, 4\", \"968\u00f73=322, 2\"),\n    @(\"467\u00f73=155, 2\", \"690\u00f77=98, 4\"),\n    @(\"203\u00f73=67, 2\", \"356\u00f77=50, 6\"),\n    @(\"289\u00f78=36, 1\", \"793\u00f77=113, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
